# Add more test data (rows 5-11) to the orders sheet, and replace the
# existing "Charl" row (row 4) with a new "Sami" order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone the existing template rows into the new rows 5-11 -------
# Row pattern repeats: Mohamed-row, Joe-row, Charl-row, Mohamed-row,
# Joe-row, Charl-row, Mohamed-row (productNum is the only thing that
# changes between repeats of the same "person" row).
$ws.Range("A2:J2").Copy($ws.Range("A5:J5"))    # Mohamed
$ws.Range("A3:J3").Copy($ws.Range("A6:J6"))    # Joe
$ws.Range("A4:J4").Copy($ws.Range("A7:J7"))    # Charl (copied before row 4 is overwritten below)
$ws.Range("A2:J2").Copy($ws.Range("A8:J8"))    # Mohamed
$ws.Range("A3:J3").Copy($ws.Range("A9:J9"))    # Joe
$ws.Range("A4:J4").Copy($ws.Range("A10:J10"))  # Charl
$ws.Range("A2:J2").Copy($ws.Range("A11:J11"))  # Mohamed

# --- 2. Fix up the productNum for every new/changed row ---------------
$ws.Range("B4").Value = 119
$ws.Range("B5").Value = 229
$ws.Range("B6").Value = 112
$ws.Range("B7").Value = 110
$ws.Range("B8").Value = 227
$ws.Range("B9").Value = 128
$ws.Range("B10").Value = 145
$ws.Range("B11").Value = 235

# --- 3. Turn row 4 into the new "Sami" order ---------------------------
$ws.Range("C4").Value = "1, 3"
$ws.Range("E4").Value = "Sami"
$ws.Range("F4").Value = "sami.basta@gmail.com"

# --- 4. Rebuild the mailto hyperlinks for every email cell -------------
# (the host's Hyperlinks.Delete() clears the whole collection, so
# everything is removed once and then re-added in document order)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:mohamed.arafa.swt@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:joe.adrian@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:sami.basta@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:mohamed.arafa.swt@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:joe.adrian@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:charl.basta@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:mohamed.arafa.swt@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:joe.adrian@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:charl.basta@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "mailto:mohamed.arafa.swt@gmail.com") | Out-Null

# --- 5. Match the final selection left by the author -------------------
$ws.Range("H11").Select() | Out-Null
